# Update average_county_temperature (column K) with new NOAA-derived data,
# and recompute the dependent worst_ashp_cop (R) / best_ashp_cop (S) values
# for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 11).Value = -3.222222222222223
$ws.Cells.Item(2, 18).Value = 1.554711451758341
$ws.Cells.Item(2, 19).Value = 1.669946025515211

# Row 3
$ws.Cells.Item(3, 11).Value = -3.222222222222223
$ws.Cells.Item(3, 18).Value = 1.475542118432027
$ws.Cells.Item(3, 19).Value = 1.575

# Row 4
$ws.Cells.Item(4, 11).Value = 0.4166666666666667
$ws.Cells.Item(4, 18).Value = 1.602020905923345
$ws.Cells.Item(4, 19).Value = 1.725399239543726

# Row 5
$ws.Cells.Item(5, 11).Value = 0.4166666666666667
$ws.Cells.Item(5, 18).Value = 1.516977491961415
$ws.Cells.Item(5, 19).Value = 1.622926829268293

# Row 6
$ws.Cells.Item(6, 11).Value = 12.51681286549706
$ws.Cells.Item(6, 18).Value = 1.782371783972741
$ws.Cells.Item(6, 19).Value = 1.939565227172176

# Row 7
$ws.Cells.Item(7, 11).Value = 12.51681286549706
$ws.Cells.Item(7, 18).Value = 1.673218141204726
$ws.Cells.Item(7, 19).Value = 1.805631235675374

# Row 12
$ws.Cells.Item(12, 11).Value = 12.51681286549706
$ws.Cells.Item(12, 18).Value = 1.782371783972741
$ws.Cells.Item(12, 19).Value = 1.939565227172176

# Row 13
$ws.Cells.Item(13, 11).Value = 12.51681286549706
$ws.Cells.Item(13, 18).Value = 1.673218141204726
$ws.Cells.Item(13, 19).Value = 1.805631235675374

# Row 16
$ws.Cells.Item(16, 11).Value = 12.51681286549706
$ws.Cells.Item(16, 18).Value = 1.782371783972741
$ws.Cells.Item(16, 19).Value = 1.939565227172176

# Row 17
$ws.Cells.Item(17, 11).Value = 12.51681286549706
$ws.Cells.Item(17, 18).Value = 1.673218141204726
$ws.Cells.Item(17, 19).Value = 1.805631235675374

# Row 18
$ws.Cells.Item(18, 11).Value = 0.4166666666666667
$ws.Cells.Item(18, 18).Value = 1.602020905923345
$ws.Cells.Item(18, 19).Value = 1.725399239543726

# Row 19
$ws.Cells.Item(19, 11).Value = 0.4166666666666667
$ws.Cells.Item(19, 18).Value = 1.516977491961415
$ws.Cells.Item(19, 19).Value = 1.622926829268293

# Row 24
$ws.Cells.Item(24, 11).Value = 15.74228395061728
$ws.Cells.Item(24, 18).Value = 1.837513876759573
$ws.Cells.Item(24, 19).Value = 2.005936573945218

# Row 25
$ws.Cells.Item(25, 11).Value = 15.74228395061728
$ws.Cells.Item(25, 18).Value = 1.720452734369724
$ws.Cells.Item(25, 19).Value = 1.861492917301914

# Row 28
$ws.Cells.Item(28, 11).Value = 15.74228395061728
$ws.Cells.Item(28, 18).Value = 1.837513876759573
$ws.Cells.Item(28, 19).Value = 2.005936573945218

# Row 29
$ws.Cells.Item(29, 11).Value = 15.74228395061728
$ws.Cells.Item(29, 18).Value = 1.861492917301914
$ws.Cells.Item(29, 19).Value = 1.861492917301914
